# Applies the "max_displacement" / "n_extra" re-ordering change described in the
# commit message:
#  - "n_extra" is moved from the bottom of the main "*" filter block up to the
#    top (directly after "t_end"), pushing n_blur..area_penalty down by one row.
#  - A brand new "max_displacement" parameter (value None) is inserted right
#    after "area_penalty" (and before "morph_transform"), pushing the
#    remaining rows (morph_transform, reversal_threshold, and the *_SF_* /
#    *dark* blocks) down by one row.
#  - While re-flowing, the "*dark*" / "bkg.secondary_factor" value changes
#    from 5 to 3, and the "*dark*" / "min_area" value changes from 20 to 25.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Right single quotation mark (U+2019), used in the "'dark..." string.
$rsquo = [char]0x2019

# Only one brand new row is actually added overall (for "max_displacement");
# "n_extra" is simply relocated within the existing 10-row block above it, so
# that block keeps the same number of rows. Insert the single new row so the
# sheet grows to the correct size (A1:D32) and everything below it shifts
# down, just like it would if a user right-clicked and chose "Insert" in the
# UI.
$ws.Rows.Item(21).Insert()

# Now (re)write every row from 11 to 30 with its final contents. Doing this
# explicitly -- rather than relying only on the shifts above -- guarantees
# the end state matches exactly regardless of how the two inserts left
# things, since several rows need their label/value pairs rewritten anyway.
$rows = @(
    @(11, "*",      "n_extra",             1),
    @(12, "*",      "n_blur",              7),
    @(13, "*",      "block_size",          31),
    @(14, "*",      "threshold_offset",    20),
    @(15, "*",      "min_area",            50),
    @(16, "*",      "max_area",            800),
    @(17, "*",      "ideal_area",          150),
    @(18, "*",      "max_aspect",          10),
    @(19, "*",      "ideal_aspect",        5),
    @(20, "*",      "area_penalty",        0.2),
    @(21, "*",      "max_displacement",    "None"),
    @(22, "*",      "morph_transform",     "[]"),
    @(23, "*",      "reversal_threshold",  0.5),
    @(24, "*_SF_*", "threshold_offset",    25),
    @(25, "*_SF_*", "min_area",            50),
    @(26, "*_SF_*", "max_aspect",          15),
    @(27, "*dark*", "threshold_offset",    10),
    @(28, "*dark*", "bkg.secondary_factor",3),
    @(29, "*dark*", "min_area",            25),
    @(30, "*dark*", "bkg.object_type",     "''dark$rsquo")
)

foreach ($row in $rows) {
    $r = $row[0]
    $ws.Cells.Item($r, 1).Value = $row[1]
    $ws.Cells.Item($r, 2).Value = $row[2]
    $ws.Cells.Item($r, 3).Value = $row[3]
}

# Row 31 is a blank placeholder row (no values), and row 32 holds the
# trailing blank-but-styled cells -- matching the original trailing blank
# rows that got pushed down by the insert above.
$ws.Range("A31:D31").ClearContents()
$ws.Range("A32:D32").ClearContents()
$ws.Range("A32:C32").Value = ""
